$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C, shifting the existing "week" column
# (old C, e.g. Jun_10) to the right to make room for two more weekly snapshots.
$ws.Range("C1:D1").EntireColumn.Insert()

# New header row: most-recent week first (B), oldest last (E).
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# Fill the two new week columns with the same "UN" (unchanged) marker used
# in column B for every ticker row, matching the existing report convention.
for ($r = 2; $r -le 27; $r++) {
    $marker = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $marker
    $ws.Cells.Item($r, 4).Value = $marker
}

Write-Host "done"
